$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 700
$ws.Range("E4").Value = 550
$ws.Range("E5").Value = 450
$ws.Range("E7").Value = 600
$ws.Range("E9").Value = 500
$ws.Range("E10").Value = 750
$ws.Range("E11").Value = 216.67
$ws.Range("E12").Value = 750
$ws.Range("E13").Value = 750
$ws.Range("E15").Value = 800
$ws.Range("E16").Value = 700
$ws.Range("E17").Value = 800
$ws.Range("E18").Value = 800
$ws.Range("E19").Value = 800
$ws.Range("E20").Value = 800
$ws.Range("E21").Value = 600
